$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.162.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.82%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.908.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.05%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.21%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'327.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.46%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.24%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4648"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.31%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3927"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.09%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'46.92"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.02%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.07962"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.12%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'1.003"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +2.74%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'22.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.81%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.927.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +3.90%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'7.139"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.01%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'5.799"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.06957"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.43%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'88.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.74%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'1.003"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.29%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.00001011"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.78%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'17.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.32%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -0.18%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'29.199.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +2.00%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.372"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.68%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'11.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.66%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.158.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +3.70%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.059"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.69%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'156.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +2.24%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'19.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.46%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'5.875"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.20%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'2.008"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.99%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'119.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.01%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.09415"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.83%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.9239"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.99%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'5.364"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.82%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.58%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'3.264"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.01%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.05854"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.29%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.167"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.20%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'8.035"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +3.45%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.02104"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.12%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.5765"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.38%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1813"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.36%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'10.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.58%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'12.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +3.57%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.5432"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.24%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.232"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +4.25%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.07097"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.29%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +2.87%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.565"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +6.72%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'112.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.60%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.077"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -5.53%  "
$ws.Range("E51").Style = "Normal"
